# chore: update Sheets via scheduled runner
# Refresh market-price driven columns (currentAveragePrice / NQ / HQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) for the leves
# whose underlying item prices moved since the last scheduled run.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 25000
$ws.Range("J13").Value = 25000
$ws.Range("L13").Value = 25000
$ws.Range("N13").Value = -25338

# Row 16
$ws.Range("H16").Value = 16203.4
$ws.Range("I16").Value = 12172.667
$ws.Range("K16").Value = 12172.667
$ws.Range("M16").Value = -11942.667

# Row 19
$ws.Range("H19").Value = 2445.5217
$ws.Range("I19").Value = 2160.2307
$ws.Range("K19").Value = 2160.2307
$ws.Range("M19").Value = -1985.2307

# Row 31
$ws.Range("H31").Value = 6745
$ws.Range("I31").Value = 6745
$ws.Range("K31").Value = 20235
$ws.Range("M31").Value = -20005

# Row 32
$ws.Range("H32").Value = 9952.546
$ws.Range("I32").Value = 11063.866
$ws.Range("J32").Value = 7571.143
$ws.Range("K32").Value = 11063.866
$ws.Range("L32").Value = 7571.143
$ws.Range("M32").Value = -10737.866
$ws.Range("N32").Value = -8223.143

# Row 76
$ws.Range("H76").Value = 8618.421
$ws.Range("J76").Value = 7100.75
$ws.Range("L76").Value = 7100.75
$ws.Range("N76").Value = -7730.75

# Row 79
$ws.Range("H79").Value = 8618.421
$ws.Range("J79").Value = 7100.75
$ws.Range("L79").Value = 7100.75
$ws.Range("N79").Value = -9284.75

# Row 88
$ws.Range("H88").Value = 3258.7
$ws.Range("I88").Value = 3319.1428
$ws.Range("J88").Value = 3226.1538
$ws.Range("K88").Value = 3319.1428
$ws.Range("L88").Value = 3226.1538
$ws.Range("M88").Value = -2913.1428
$ws.Range("N88").Value = -4038.1538

# Row 91
$ws.Range("H91").Value = 3258.7
$ws.Range("I91").Value = 3319.1428
$ws.Range("J91").Value = 3226.1538
$ws.Range("K91").Value = 3319.1428
$ws.Range("L91").Value = 3226.1538
$ws.Range("M91").Value = -1915.1428
$ws.Range("N91").Value = -6034.1538

# Row 137
$ws.Range("H137").Value = 1469.84
$ws.Range("I137").Value = 996.0714
$ws.Range("J137").Value = 2072.818
$ws.Range("K137").Value = 2988.2142
$ws.Range("L137").Value = 6218.454000000001
$ws.Range("M137").Value = -438.2142000000003
$ws.Range("N137").Value = -11318.454


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2772.8723
$ws.Range("I32").Value = 775.06976
$ws.Range("J32").Value = 24249.25
$ws.Range("K32").Value = 775.06976
$ws.Range("L32").Value = 24249.25
$ws.Range("M32").Value = -488.06976
$ws.Range("N32").Value = -24823.25

# Row 63
$ws.Range("H63").Value = 1643.125
$ws.Range("I63").Value = 1856.7142
$ws.Range("J63").Value = 148
$ws.Range("K63").Value = 1856.7142
$ws.Range("L63").Value = 148
$ws.Range("M63").Value = -1170.7142
$ws.Range("N63").Value = -1520

# Row 66
$ws.Range("H66").Value = 1643.125
$ws.Range("I66").Value = 1856.7142
$ws.Range("J66").Value = 148
$ws.Range("K66").Value = 9283.571
$ws.Range("L66").Value = 740
$ws.Range("M66").Value = -5851.571
$ws.Range("N66").Value = -7604

# Row 88
$ws.Range("H88").Value = 2080.8333
$ws.Range("I88").Value = 1673.5
$ws.Range("J88").Value = 2284.5
$ws.Range("K88").Value = 1673.5
$ws.Range("L88").Value = 2284.5
$ws.Range("M88").Value = -1267.5
$ws.Range("N88").Value = -3096.5

# Row 91
$ws.Range("H91").Value = 2080.8333
$ws.Range("I91").Value = 1673.5
$ws.Range("J91").Value = 2284.5
$ws.Range("K91").Value = 1673.5
$ws.Range("L91").Value = 2284.5
$ws.Range("M91").Value = -269.5
$ws.Range("N91").Value = -5092.5

# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0

# Row 122
$ws.Range("H122").Value = 3206.2856
$ws.Range("I122").Value = 3206.2856
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9618.856800000001
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -7168.856800000001

# Row 132
$ws.Range("H132").Value = 3510.4
$ws.Range("I132").Value = 3888.375
$ws.Range("K132").Value = 11665.125
$ws.Range("M132").Value = -9135.125


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2341.6924
$ws.Range("I86").Value = 1692.0952
$ws.Range("K86").Value = 1692.0952
$ws.Range("M86").Value = -569.0952

# Row 89
$ws.Range("H89").Value = 2341.6924
$ws.Range("I89").Value = 1692.0952
$ws.Range("K89").Value = 8460.476000000001
$ws.Range("M89").Value = -2844.476000000001

# Row 96
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 233976.08
$ws.Range("I122").Value = 252890.75
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 758672.25
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -756222.25
$ws.Range("N122").Value = -25900


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 1522.1875
$ws.Range("I75").Value = 2097.5
$ws.Range("J75").Value = 1440
$ws.Range("K75").Value = 6292.5
$ws.Range("L75").Value = 4320
$ws.Range("M75").Value = -5294.5
$ws.Range("N75").Value = -6316

# Row 78
$ws.Range("H78").Value = 1522.1875
$ws.Range("I78").Value = 2097.5
$ws.Range("J78").Value = 1440
$ws.Range("K78").Value = 18877.5
$ws.Range("L78").Value = 12960
$ws.Range("M78").Value = -13885.5
$ws.Range("N78").Value = -22944

# Row 114
$ws.Range("H114").Value = 1965.4
$ws.Range("I114").Value = 1999
$ws.Range("J114").Value = 1957
$ws.Range("K114").Value = 5997
$ws.Range("L114").Value = 5871
$ws.Range("M114").Value = -2743
$ws.Range("N114").Value = -12379

# Row 131
$ws.Range("H131").Value = 2401.3914
$ws.Range("I131").Value = 1358
$ws.Range("J131").Value = 3357.8333
$ws.Range("K131").Value = 4074
$ws.Range("L131").Value = 10073.4999
$ws.Range("M131").Value = 966
$ws.Range("N131").Value = -20153.4999


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 116
$ws.Range("H116").Value = 89900
$ws.Range("J116").Value = 89900
$ws.Range("L116").Value = 89900
$ws.Range("N116").Value = -99078

# Row 122
$ws.Range("H122").Value = 1102.4
$ws.Range("J122").Value = 951
$ws.Range("L122").Value = 2853
$ws.Range("N122").Value = -7753


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3900.5
$ws.Range("J68").Value = 8335.333000000001
$ws.Range("L68").Value = 8335.333000000001
$ws.Range("N68").Value = -9833.333000000001

# Row 71
$ws.Range("H71").Value = 3900.5
$ws.Range("J71").Value = 8335.333000000001
$ws.Range("L71").Value = 41676.665
$ws.Range("N71").Value = -49164.665

# Row 122
$ws.Range("H122").Value = 14231.75
$ws.Range("I122").Value = 14231.75
$ws.Range("K122").Value = 42695.25
$ws.Range("M122").Value = -40245.25


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 9
$ws.Range("H9").Value = 3028.8572
$ws.Range("I9").Value = 1040.4
$ws.Range("J9").Value = 8000
$ws.Range("K9").Value = 1040.4
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = -900.4000000000001
$ws.Range("N9").Value = -8280

# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

# Row 122
$ws.Range("H122").Value = 2217.5
$ws.Range("I122").Value = 2127.8333
$ws.Range("K122").Value = 6383.499899999999
$ws.Range("M122").Value = -3933.499899999999

